$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) --------------------------------------------------
# Bold, centered (horizontal + top vertical), thin-boxed header style.
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
    $cell.Borders.Weight = 2            # xlThin
}

# --- Data row (row 2) -----------------------------------------------------
$ws.Cells.Item(2, 1).Value = "JD_001"
$ws.Cells.Item(2, 2).Value = "Junior RPA Developer"
$ws.Cells.Item(2, 3).Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = "Hybrid"
$ws.Cells.Item(2, 7).Value = "Pune, Maharashtra, India"

# The embedded line break in C2 triggers an auto row-height bump;
# restore the row to its default auto-fit height.
$ws.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# Leave the selection on A1 (matches a freshly populated sheet).
$ws.Range("A1").Select() | Out-Null
